$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.642.86'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '3.710.10'
$ws.Range('E3').Value = '  +1.03%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = "'673.05"
$ws.Range('E5').Value = '  -1.46%  '
$ws.Range('D6').Value = "'161.92"
$ws.Range('E6').Value = '  +2.69%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('E9').Value = '  +1.04%  '
$ws.Range('D10').Value = "'7.11"
$ws.Range('E10').Value = '  +2.11%  '
$ws.Range('D11').Value = "'0.444"
$ws.Range('E11').Value = '  +2.18%  '
$ws.Range('E12').Value = '  +1.55%  '
$ws.Range('D13').Value = "'32.89"
$ws.Range('E13').Value = '  +2.44%  '
$ws.Range('D14').Value = '3.724.88'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').Value = '69.681.94'
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('E16').Value = '  +1.64%  '
$ws.Range('D17').Value = "'16.31"
$ws.Range('E17').Value = '  +3.32%  '
$ws.Range('D18').Value = "'6.51"
$ws.Range('E18').Value = '  +2.33%  '
$ws.Range('D19').Value = "'474.06"
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('D20').Value = "'9.81"
$ws.Range('E20').Value = '  -1.10%  '
$ws.Range('D21').Value = "'0.654"
$ws.Range('E21').Value = '  +1.25%  '
$ws.Range('D22').Value = "'80.49"
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('D23').Value = '3.861.27'
$ws.Range('E23').Value = '  +1.11%  '
$ws.Range('E24').Value = '  +5.50%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').Value = "'11.02"
$ws.Range('E26').Value = '  +1.09%  '
$ws.Range('E27').Value = '  +0.51%  '
$ws.Range('D28').Value = "'2.69"
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('E30').Value = '  +2.22%  '
$ws.Range('E31').Value = '  +7.42%  '
$ws.Range('E32').Value = '  +1.23%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').Value = "'26.94"
$ws.Range('E34').Value = '  +0.45%  '
$ws.Range('D35').Value = '3.700.58'
$ws.Range('E35').Value = '  +1.40%  '
$ws.Range('D36').Value = "'8.57"
$ws.Range('E36').Value = '  +5.26%  '
$ws.Range('D37').Value = "'6.11"
$ws.Range('E37').Value = '  +0.89%  '
$ws.Range('D39').Value = "'2.26"
$ws.Range('E39').Value = '  +2.75%  '
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('D41').Value = "'0.0915"
$ws.Range('E41').Value = '  +2.19%  '
$ws.Range('D42').Value = "'174.07"
$ws.Range('D43').Value = "'0.941"
$ws.Range('E43').Value = '  +0.26%  '
$ws.Range('D44').Value = "'47.11"
$ws.Range('E44').Value = '  -0.91%  '
$ws.Range('D45').Value = "'2.78"
$ws.Range('E45').Value = '  +3.32%  '
$ws.Range('D46').Value = "'0.000282"
$ws.Range('E46').Value = '  +1.58%  '
$ws.Range('E47').Value = '  +2.28%  '
$ws.Range('D48').Value = "'27.72"
$ws.Range('E48').Value = '  +3.60%  '
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('E50').Value = '  +1.93%  '
$ws.Range('D51').Value = "'0.268"
$ws.Range('E51').Value = '  +1.87%  '
